$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 with the same text and formatting as the other
# header cells (A1:C1), which carry bold font, border and centered alignment.
$ws.Range("D1").Value = "Tipo"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Update the recalculated MSE / R2 values in row 2
$ws.Range("B2").Value = 0.1256512455282837
$ws.Range("C2").Value = 0.9907164282362919

# Add the new "Tipo" value for row 2 (no special formatting, like A2/B2/C2)
$ws.Range("D2").Value = "single"

# Refresh the used range so the sheet dimension reflects the new column D
$ws.UsedRange | Out-Null
